{"js": "// The second paragraph of the document contains a Word field\n// ( { m:'doc.html'.fromHTMLURI() } ) encoded as fldChar/instrText runs\n// (begin, instrText tokens, end) around a \"_GoBack\" bookmark. The edit\n// turns that field code into plain visible text runs that spell out the\n// same token stream, literally showing \"{\", \"m\", \":\", \"'\", \"doc.html\",\n// (bookmark kept in place), \"'.fromHTMLURI()\" and \"}\" - i.e. the field\n// delimiters/instrText are replaced by w:t runs with equivalent text,\n// and the stray leading/trailing space tokens are dropped, being folded\n// into the \"{\" and \"}\" runs.\n//\n// Office.js has no direct way to manipulate w:fldChar / w:instrText runs,\n// so we rebuild the whole paragraph via OOXML (keeping the bookmark) and\n// swap it in with insertOoxml/replace.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Find the paragraph that holds the field (it is the one whose\n// underlying XML contains the \"doc.html\" instruction text / fldChar).\nconst ooxmlResults = [];\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  ooxmlResults.push(paragraphs.items[i].getOoxml());\n}\nawait context.sync();\n\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const xml = ooxmlResults[i].value;\n  if (xml && xml.indexOf(\"doc.html\") !== -1 && xml.indexOf(\"fldChar\") !== -1) {\n    target = paragraphs.items[i];\n    break;\n  }\n}\n\nif (target) {\n  const paragraphXml =\n    '<w:p xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\" ' +\n      'w:rsidR=\"00C52979\" w:rsidRDefault=\"00C52979\" w:rsidP=\"00F5495F\">' +\n      \"<w:r><w:t>{</w:t></w:r>\" +\n      \"<w:r><w:t>m</w:t></w:r>\" +\n      \"<w:r><w:t>:</w:t></w:r>\" +\n      \"<w:r><w:t>'</w:t></w:r>\" +\n      \"<w:r><w:t>doc.html</w:t></w:r>\" +\n      '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>' +\n      '<w:bookmarkEnd w:id=\"0\"/>' +\n      \"<w:r><w:t>'.fromHTMLURI()</w:t></w:r>\" +\n      '<w:r><w:t xml:space=\"preserve\">}</w:t></w:r>' +\n    \"</w:p>\";\n\n  const packageXml =\n    '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n      '<pkg:part pkg:name=\"/word/document.xml\" ' +\n        'pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n        \"<pkg:xmlData>\" +\n          '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n            \"<w:body>\" + paragraphXml + \"</w:body>\" +\n          \"</w:document>\" +\n        \"</pkg:xmlData>\" +\n      \"</pkg:part>\" +\n    \"</pkg:package>\";\n\n  target.insertOoxml(packageXml, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# The document contains one Word field: { m:'doc.html'.fromHTMLURI() }\n# written as fldChar(begin)/instrText.../fldChar(end) runs wrapped around a\n# \"_GoBack\" bookmark, inside its own paragraph. The edit turns that field\n# code into plain, always-visible text runs that spell out the same token\n# stream (\"{\", \"m\", \":\", \"'\", \"doc.html\", the bookmark, \"'.fromHTMLURI()\",\n# \"}\"), i.e. the field is removed and replaced by literal w:t runs, with\n# the stray leading/trailing single-space instrText tokens folded into the\n# \"{\" and \"}\" runs.\n\n$d = $word.ActiveDocument\n\n# Locate the field whose instruction text references doc.html.\n$targetField = $null\nforeach ($f in $d.Fields) {\n    if ($f.Code.Text -like \"*doc.html*\") {\n        $targetField = $f\n        break\n    }\n}\n\nif ($targetField -ne $null) {\n    $fieldStart = $targetField.Code.Start\n\n    # Find the paragraph that owns this field.\n    $targetParagraph = $null\n    foreach ($p in $d.Paragraphs) {\n        if ($p.Range.Start -le $fieldStart -and $fieldStart -lt $p.Range.End) {\n            $targetParagraph = $p\n            break\n        }\n    }\n\n    if ($targetParagraph -ne $null) {\n        $insertionRange = $targetParagraph.Range\n\n        # Remove the field (fldChar begin/end + instrText runs + the\n        # bookmark that lived inside it) entirely before rebuilding the\n        # paragraph content.\n        $targetField.Delete()\n\n        # Rebuild the paragraph as plain text runs, keeping the\n        # \"_GoBack\" bookmark in the same relative spot, via raw OOXML so\n        # that each literal token becomes its own w:t run, matching the\n        # original instrText run boundaries.\n        $paragraphXml =\n            '<w:p xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\" ' +\n              'w:rsidR=\"00C52979\" w:rsidRDefault=\"00C52979\" w:rsidP=\"00F5495F\">' +\n              '<w:r><w:t>{</w:t></w:r>' +\n              '<w:r><w:t>m</w:t></w:r>' +\n              '<w:r><w:t>:</w:t></w:r>' +\n              \"<w:r><w:t>'</w:t></w:r>\" +\n              '<w:r><w:t>doc.html</w:t></w:r>' +\n              '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>' +\n              '<w:bookmarkEnd w:id=\"0\"/>' +\n              \"<w:r><w:t>'.fromHTMLURI()</w:t></w:r>\" +\n              '<w:r><w:t xml:space=\"preserve\">}</w:t></w:r>' +\n            '</w:p>'\n\n        $packageXml =\n            '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n            '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n              '<pkg:part pkg:name=\"/word/document.xml\" ' +\n                'pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n                '<pkg:xmlData>' +\n                  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n                    '<w:body>' + $paragraphXml + '</w:body>' +\n                  '</w:document>' +\n                '</pkg:xmlData>' +\n              '</pkg:part>' +\n            '</pkg:package>'\n\n        $insertionRange.InsertXML($packageXml)\n    }\n}\n"}
